$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Select()

$f3 = '=IF(OR(A3<$H$2,A3>$I$2),"OUTLIER","NORMAL")'
$ws.Range("B3").Formula = $f3

$f4 = '=IF(OR(A4<$H$2,A4>$I$2),"OUTLIER","NORMAL")'
$ws.Range("B4").Formula = $f4

$ws.Range("D4").Value = " "
$ws.Range("F9").Value = " "
